# "The updated SIQ Reviewing"
# Update the Status column (J) for the first four review rows from their
# previous values ("Approved" / "Rejected") to "Pending", and move the
# active selection/scroll position to reflect where the reviewer was
# working (around J10, scrolled so B2 is the top-left visible cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RVW-001 .. RVW-004 status -> Pending
$ws.Range("J2").Value = "Pending"
$ws.Range("J3").Value = "Pending"
$ws.Range("J4").Value = "Pending"
$ws.Range("J5").Value = "Pending"

# Scroll the view so B2 is the top-left visible cell (best-effort; some
# hosts don't persist plain scroll position without freeze panes) and move
# the selection to J10, matching the reviewer's final cursor position.
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
    # Not fatal if the host doesn't support direct scroll positioning.
}

$ws.Range("J10").Select()
